$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 50")

# Row 13: new logboek entry "Debug"
$ws.Range("C13").Value = 0.51388888888888895
$ws.Range("D13").Value = 0.53125
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = "Debug"

# Row 14: new logboek entry describing the work done
$ws.Range("C14").Value = 0.5625
$ws.Range("D14").Value = 0.61458333333333337
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = "Knoppen werkend gemaakt met edge detection. Begin gemaakt aan Scene's dat als je enter drukt naar de juiste scene gaat."

# The long activity text wraps onto multiple lines, so the row grows taller
$ws.Rows.Item(14).RowHeight = 42.75

# Move the active selection to F14 (where the user last edited)
$ws.Range("F14").Select()
